# Repull data, push all data, mean calculation
# Update the dSF column (F) values for rows 3, 4, 5, 6, 8, 14

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = -2
$ws.Range("F4").Value = -2
$ws.Range("F5").Value = -7
$ws.Range("F6").Value = -4
$ws.Range("F8").Value = 7
$ws.Range("F14").Value = 3
